$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Copy style from an existing header cell (F1) to the new header cells
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Fill data rows 2-5 for columns G (Elapsed Time) and H (CPU)
for ($r = 2; $r -le 5; $r++) {
    $ws.Cells.Item($r, 7).Value = 0.4794827245333484
    $ws.Cells.Item($r, 8).Value = 0.996
}
